# Scheduled runner update: refresh market-price/profit columns (H:N) on
# the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1825.9286
$ws.Range("J40").Value = 1900.3334
$ws.Range("L40").Value = 1900.3334
$ws.Range("N40").Value = -2250.3334

$ws.Range("H113").Value = 2457.1304
$ws.Range("I113").Value = 2581
$ws.Range("J113").Value = 2224.875
$ws.Range("K113").Value = 2581
$ws.Range("L113").Value = 2224.875
$ws.Range("M113").Value = 673
$ws.Range("N113").Value = -8732.875

$ws.Range("H129").Value = 850
$ws.Range("I129").Value = 253.3
$ws.Range("J129").Value = 1276.2142
$ws.Range("K129").Value = 759.9000000000001
$ws.Range("L129").Value = 3828.6426
$ws.Range("M129").Value = 4240.1
$ws.Range("N129").Value = -13828.6426

$ws.Range("H132").Value = 21942.389
$ws.Range("I132").Value = 23983.568
$ws.Range("K132").Value = 71950.704
$ws.Range("M132").Value = -69420.704

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2112.625
$ws.Range("I45").Value = 2043.3939
$ws.Range("J45").Value = 2264.9333
$ws.Range("K45").Value = 2043.3939
$ws.Range("L45").Value = 2264.9333
$ws.Range("M45").Value = -1666.3939
$ws.Range("N45").Value = -3018.9333

$ws.Range("H132").Value = 27807352
$ws.Range("I132").Value = 38463030
$ws.Range("J132").Value = 102584.2
$ws.Range("K132").Value = 115389090
$ws.Range("L132").Value = 307752.6
$ws.Range("M132").Value = -115386560
$ws.Range("N132").Value = -312812.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 911.14703
$ws.Range("I107").Value = 600.375
$ws.Range("K107").Value = 600.375
$ws.Range("M107").Value = 1319.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 5909.9
$ws.Range("I55").Value = 2374.75
$ws.Range("J55").Value = 8266.666999999999
$ws.Range("K55").Value = 2374.75
$ws.Range("L55").Value = 8266.666999999999
$ws.Range("M55").Value = -2059.75
$ws.Range("N55").Value = -8896.666999999999

$ws.Range("H62").Value = 3780
$ws.Range("I62").Value = 2225
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 2225
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -1601
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 3780
$ws.Range("I65").Value = 2225
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 11125
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -8005
$ws.Range("N65").Value = -56240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1010.05884
$ws.Range("I5").Value = 979.4375
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 2938.3125
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -2826.3125
$ws.Range("N5").Value = -4724

$ws.Range("H122").Value = 788
$ws.Range("I122").Value = 330.33334
$ws.Range("J122").Value = 1474.5
$ws.Range("K122").Value = 2973.00006
$ws.Range("L122").Value = 13270.5
$ws.Range("M122").Value = -523.0000600000003
$ws.Range("N122").Value = -18170.5

$ws.Range("H131").Value = 5748116
$ws.Range("I131").Value = 2800
$ws.Range("J131").Value = 5883299.5
$ws.Range("K131").Value = 8400
$ws.Range("L131").Value = 17649898.5
$ws.Range("M131").Value = -3360
$ws.Range("N131").Value = -17659978.5

$ws.Range("H135").Value = 1010.05884
$ws.Range("I135").Value = 979.4375
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 8814.9375
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -6279.9375
$ws.Range("N135").Value = -18570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1406.8182
$ws.Range("I113").Value = 1296.6666
$ws.Range("J113").Value = 1642.8572
$ws.Range("K113").Value = 1296.6666
$ws.Range("L113").Value = 1642.8572
$ws.Range("M113").Value = 873.3334
$ws.Range("N113").Value = -5982.8572

$ws.Range("H132").Value = 32272.908
$ws.Range("I132").Value = 1372.0526
$ws.Range("J132").Value = 74209.78999999999
$ws.Range("K132").Value = 4116.1578
$ws.Range("L132").Value = 222629.37
$ws.Range("M132").Value = -1586.1578
$ws.Range("N132").Value = -227689.37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 709.3333
$ws.Range("I22").Value = 491.79166
$ws.Range("J22").Value = 957.9524
$ws.Range("K22").Value = 491.79166
$ws.Range("L22").Value = 957.9524
$ws.Range("M22").Value = -196.79166
$ws.Range("N22").Value = -1547.9524

$ws.Range("H27").Value = 709.3333
$ws.Range("I27").Value = 491.79166
$ws.Range("J27").Value = 957.9524
$ws.Range("K27").Value = 491.79166
$ws.Range("L27").Value = 957.9524
$ws.Range("M27").Value = -384.79166
$ws.Range("N27").Value = -1171.9524

$ws.Range("H82").Value = 2146.4614
$ws.Range("I82").Value = 2583.3333
$ws.Range("J82").Value = 1772
$ws.Range("K82").Value = 2583.3333
$ws.Range("L82").Value = 1772
$ws.Range("M82").Value = -2222.3333
$ws.Range("N82").Value = -2494

$ws.Range("H85").Value = 2146.4614
$ws.Range("I85").Value = 2583.3333
$ws.Range("J85").Value = 1772
$ws.Range("K85").Value = 2583.3333
$ws.Range("L85").Value = 1772
$ws.Range("M85").Value = -1335.3333
$ws.Range("N85").Value = -4268

$ws.Range("H136").Value = 17844
$ws.Range("I136").Value = 18484.666
$ws.Range("J136").Value = 16562.666
$ws.Range("K136").Value = 55453.99800000001
$ws.Range("L136").Value = 49687.99800000001
$ws.Range("M136").Value = -52903.99800000001
$ws.Range("N136").Value = -54787.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1221.4667
$ws.Range("I122").Value = 1213.6316
$ws.Range("J122").Value = 1235
$ws.Range("K122").Value = 3640.8948
$ws.Range("L122").Value = 3705
$ws.Range("M122").Value = -1190.8948
$ws.Range("N122").Value = -8605

$ws.Range("H132").Value = 79045220
$ws.Range("I132").Value = 125556640
$ws.Range("J132").Value = 2935624.2
$ws.Range("K132").Value = 376669920
$ws.Range("L132").Value = 8806872.600000001
$ws.Range("M132").Value = -376667390
$ws.Range("N132").Value = -8811932.600000001

$ws.Range("H136").Value = 25955.875
$ws.Range("I136").Value = 50836.25
$ws.Range("J136").Value = 1075.5
$ws.Range("K136").Value = 152508.75
$ws.Range("L136").Value = 3226.5
$ws.Range("N136").Value = -8326.5
